$wb = $excel.ActiveWorkbook

# --- Rename sheets ---
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws1.Name = "1碑影迷踪"
$ws2.Name = "2消失的龙"

# --- Update the header text (rewrites the shared strings "网"->"网址" and "提示"->"提示词" in place) ---
$ws1.Range("A1").Value = "网址"
$ws1.Range("B1").Value = "提示词"

# --- Populate sheet 2 header row with the same two strings ---
$ws2.Range("A1").Value = "网址"
$ws2.Range("B1").Value = "提示词"

# --- Selections on sheet 1 ---
[void]$ws1.Range("B2").Select()

# --- Activate sheet 2 and set its selection ---
$ws2.Activate()
[void]$ws2.Range("C8").Select()
